$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the "att" column (C2:C33) from 1 to -1
$ws.Range("C2:C33").Value = -1

# Update the view state: scroll up and move the active selection to G31
[void]$ws.Range("G31").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
